# Update the "想去人数" (want-to-go count) figures for both the "展览"
# sheet and the "全部类型" sheet, which duplicate the same rows.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 665
    $ws.Range("F7").Value = 43
    $ws.Range("F8").Value = 2618
    $ws.Range("F9").Value = 4150
}
